# Update NATMI TPM-derived LR-pair metrics (Glg1-Sele) with newly computed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 21.27021466666666
$ws.Range("H2").Value = 63.810644
$ws.Range("I2").Value = 0.1783704271312809
$ws.Range("J2").Value = 0.1783704271312809
$ws.Range("M2").Value = 12.67919733333333
$ws.Range("N2").Value = 38.037592
$ws.Range("O2").Value = 0.9871416146107245
$ws.Range("P2").Value = 0.9871416146107247
$ws.Range("Q2").Value = 269.6892490810275
$ws.Range("R2").Value = 2427.203241729248
$ws.Range("S2").Value = 0.1760768714371772
$ws.Range("T2").Value = 0.1760768714371772
$ws.Range("G3").Value = 21.27021466666666
$ws.Range("H3").Value = 63.810644
$ws.Range("I3").Value = 0.1783704271312809
$ws.Range("J3").Value = 0.1783704271312809
$ws.Range("M3").Value = 0.1651576666666667
$ws.Range("N3").Value = 0.495473
$ws.Range("O3").Value = 0.01285838538927542
$ws.Range("P3").Value = 0.01285838538927542
$ws.Range("Q3").Value = 3.512939023845778
$ws.Range("R3").Value = 31.616451214612
$ws.Range("S3").Value = 0.002293555694103678
$ws.Range("T3").Value = 0.002293555694103677
$ws.Range("G4").Value = 53.56207000000001
$ws.Range("H4").Value = 160.68621
$ws.Range("I4").Value = 0.4491675074115645
$ws.Range("J4").Value = 0.4491675074115645
$ws.Range("M4").Value = 12.67919733333333
$ws.Range("N4").Value = 38.037592
$ws.Range("O4").Value = 0.9871416146107245
$ws.Range("P4").Value = 0.9871416146107247
$ws.Range("Q4").Value = 679.1240551118135
$ws.Range("R4").Value = 6112.116496006322
$ws.Range("S4").Value = 0.4433919384969264
$ws.Range("T4").Value = 0.4433919384969264
$ws.Range("G5").Value = 53.56207000000001
$ws.Range("H5").Value = 160.68621
$ws.Range("I5").Value = 0.4491675074115645
$ws.Range("J5").Value = 0.4491675074115645
$ws.Range("M5").Value = 0.1651576666666667
$ws.Range("N5").Value = 0.495473
$ws.Range("O5").Value = 0.01285838538927542
$ws.Range("P5").Value = 0.01285838538927542
$ws.Range("Q5").Value = 8.846186503036668
$ws.Range("R5").Value = 79.61567852733
$ws.Range("S5").Value = 0.005775568914638119
$ws.Range("T5").Value = 0.005775568914638119
$ws.Range("G6").Value = 25.97301533333334
$ws.Range("H6").Value = 77.91904600000001
$ws.Range("I6").Value = 0.2178077612989132
$ws.Range("J6").Value = 0.2178077612989132
$ws.Range("M6").Value = 12.67919733333333
$ws.Range("N6").Value = 38.037592
$ws.Range("O6").Value = 0.9871416146107245
$ws.Range("P6").Value = 0.9871416146107247
$ws.Range("Q6").Value = 329.3169867530258
$ws.Range("R6").Value = 2963.852880777233
$ws.Range("S6").Value = 0.2150071051633564
$ws.Range("T6").Value = 0.2150071051633564
$ws.Range("G7").Value = 25.97301533333334
$ws.Range("H7").Value = 77.91904600000001
$ws.Range("I7").Value = 0.2178077612989132
$ws.Range("J7").Value = 0.2178077612989132
$ws.Range("M7").Value = 0.1651576666666667
$ws.Range("N7").Value = 0.495473
$ws.Range("O7").Value = 0.01285838538927542
$ws.Range("P7").Value = 0.01285838538927542
$ws.Range("Q7").Value = 4.289642608750889
$ws.Range("R7").Value = 38.606783478758
$ws.Range("S7").Value = 0.002800656135556733
$ws.Range("T7").Value = 0.002800656135556733
$ws.Range("G8").Value = 18.442128
$ws.Range("H8").Value = 55.326384
$ws.Range("I8").Value = 0.1546543041582415
$ws.Range("J8").Value = 0.1546543041582415
$ws.Range("M8").Value = 12.67919733333333
$ws.Range("N8").Value = 38.037592
$ws.Range("O8").Value = 0.9871416146107245
$ws.Range("P8").Value = 0.9871416146107247
$ws.Range("Q8").Value = 233.831380158592
$ws.Range("R8").Value = 2104.482421427329
$ws.Range("S8").Value = 0.1526656995132646
$ws.Range("T8").Value = 0.1526656995132646
$ws.Range("G9").Value = 18.442128
$ws.Range("H9").Value = 55.326384
$ws.Range("I9").Value = 0.1546543041582415
$ws.Range("J9").Value = 0.1546543041582415
$ws.Range("M9").Value = 0.1651576666666667
$ws.Range("N9").Value = 0.495473
$ws.Range("O9").Value = 0.01285838538927542
$ws.Range("P9").Value = 0.01285838538927542
$ws.Range("Q9").Value = 3.045858828848
$ws.Range("R9").Value = 27.412729459632
$ws.Range("S9").Value = 0.001988604644976888
$ws.Range("T9").Value = 0.001988604644976888